# Revert "Update Summary and BOM"
# This reverts commit bc97cbde53706c8a9a4818ee2694dbdf19bd6d33, which had
# turned the placeholder text in cell G5 of the "BOM" sheet into a live
# hyperlink pointing at the project's GitHub repo. Reverting means: drop the
# hyperlink, put the plain placeholder text back, and give the cell a
# bold/underlined "hyperlink-style" look (its text is not an actual link
# any more, but it keeps the same visual styling family).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

$cell = $ws.Range("G5")

# Remove the hyperlink that lived on G5 (also drops its relationship).
if ($cell.Hyperlinks.Count -gt 0) {
    $cell.Hyperlinks.Delete()
}

# Restore the old placeholder text.
$cell.Value = "Add github link"

# Re-apply a bold + underlined hyperlink-colored look to the cell (based on
# the built-in "Hyperlink" cell style, but bolded/underlined on top of it).
$cell.Style = "Hyperlink"
$cell.Font.Bold = $true
$cell.Font.Underline = $true

# Restore the sheet selection/scroll position that shipped with the older
# revision of the workbook.
$ws.Activate()
$ws.Range("C35").Select()
